$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four stale demand rows (2-5); row 6 shifts up to become row 2,
# carrying forward its already-empty Motif/Date_Traitement/Traite_Par/Commentaires cells.
$ws.Range("A2:I5").EntireRow.Delete()

# Update (now) row 2 with the new demand entry
$ws.Range("A2").Value = "20250526_133307"
$ws.Range("B2").Value = "2025-05-26 13:33:07"
$ws.Range("C2").Value = "Marie Martin"
$ws.Range("D2").Value = "{'chantier': 'Aluminium - Table Aluminium 02', 'urgence': 'Normal', 'date_souhaitee': '2025-05-26', 'produits': {'2140736376': {'produit': 'Crémone OB F8 mm P220  600/900 20093009', 'quantite': 4, 'emplacement': 'E1'}}}"
$ws.Range("F2").Value = "En attente"
